$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6662
$ws.Range("J3").Value = 7041
$ws.Range("J4").Value = 1529
$ws.Range("I5").Value = 625
$ws.Range("J5").Value = 552
$ws.Range("J6").Value = 9364
$ws.Range("I7").Value = 22661
$ws.Range("J7").Value = 25148

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 237
$ws.Range("J7").Value = 358

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 472
$ws.Range("J4").Value = 83
$ws.Range("J5").Value = 41
$ws.Range("J6").Value = 552
$ws.Range("J7").Value = 1579

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 265
$ws.Range("J5").Value = 47
$ws.Range("J6").Value = 397
$ws.Range("J7").Value = 1137

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 227
$ws.Range("J3").Value = 262
$ws.Range("J6").Value = 228
$ws.Range("J7").Value = 775

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 153
$ws.Range("J6").Value = 97
$ws.Range("J7").Value = 384

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J6").Value = 191
$ws.Range("J7").Value = 730
$ws.Range("J8").Value = 1579
$ws.Range("J10").Value = 184
$ws.Range("J15").Value = 297
$ws.Range("J19").Value = 735
$ws.Range("J20").Value = 524
$ws.Range("J23").Value = 231
$ws.Range("J29").Value = 1372
$ws.Range("J33").Value = 1137
$ws.Range("J36").Value = 340
$ws.Range("J37").Value = 775
$ws.Range("J41").Value = 174
$ws.Range("J42").Value = 1087
$ws.Range("J43").Value = 215
$ws.Range("J44").Value = 192
$ws.Range("J49").Value = 161
$ws.Range("J53").Value = 358
$ws.Range("J54").Value = 481
$ws.Range("J55").Value = 384
$ws.Range("J57").Value = 116
$ws.Range("J60").Value = 145
$ws.Range("I63").Value = 172
$ws.Range("J63").Value = 80
$ws.Range("J64").Value = 168
$ws.Range("J67").Value = 947
$ws.Range("J72").Value = 97
$ws.Range("J78").Value = 297
$ws.Range("J79").Value = 708
$ws.Range("J84").Value = 207
$ws.Range("J85").Value = 1046
$ws.Range("J88").Value = 260
$ws.Range("J89").Value = 322
$ws.Range("J90").Value = 268
$ws.Range("J91").Value = 291
$ws.Range("J92").Value = 79
$ws.Range("J94").Value = 269
$ws.Range("J98").Value = 185
$ws.Range("J99").Value = 384
$ws.Range("I101").Value = 22661
$ws.Range("J101").Value = 25148

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 352
$ws.Range("J4").Value = 64
$ws.Range("J7").Value = 947

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J3").Value = 69
$ws.Range("J7").Value = 207

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 161

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 117
$ws.Range("J3").Value = 101
$ws.Range("J7").Value = 481

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J6").Value = 348
$ws.Range("J7").Value = 1372

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 178
$ws.Range("J6").Value = 284
$ws.Range("J7").Value = 735

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 192

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 76
$ws.Range("J6").Value = 201

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J5").Value = 7
$ws.Range("J7").Value = 191

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J3").Value = 26
$ws.Range("J7").Value = 174

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 232
$ws.Range("J3").Value = 212
$ws.Range("J6").Value = 579
$ws.Range("J7").Value = 1087

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 184

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J5").Value = 4
$ws.Range("J7").Value = 297

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J4").Value = 17
$ws.Range("J7").Value = 384

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J2").Value = 62
$ws.Range("J7").Value = 231

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 122
$ws.Range("J7").Value = 291

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 241
$ws.Range("J7").Value = 708

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 168

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 177
$ws.Range("J7").Value = 524

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 110
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 340

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 226
$ws.Range("J7").Value = 730

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 52
$ws.Range("J3").Value = 53
$ws.Range("J7").Value = 269

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J6").Value = 129
$ws.Range("J7").Value = 297

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J3").Value = 26
$ws.Range("J6").Value = 116
$ws.Range("J7").Value = 185

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J2").Value = 55
$ws.Range("J7").Value = 260

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 93
$ws.Range("J7").Value = 322

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J3").Value = 75
$ws.Range("J7").Value = 268

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J2").Value = 29
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J3").Value = 42
$ws.Range("J6").Value = 128
$ws.Range("J7").Value = 215

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 282
$ws.Range("J6").Value = 304
$ws.Range("J7").Value = 1046

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 97
